# Auto update Excel log
#
# Appends newly-logged sensor readings to the end of the "PIR" and
# "Humidity" sheets (rows keep accumulating as the sensors report in).
#
# Date-like text (column A, e.g. "2026-01-30") and percentage text
# (column E on the Humidity sheet, e.g. "86.4%") are written with an
# explicit Text number format ("@") applied first - otherwise the COM
# layer auto-coerces those strings into date/number serials, and the
# log is meant to keep every field as plain text, same as the rest of
# the sheet.

$wb = $excel.ActiveWorkbook

# --- PIR sheet: new "Bathroom" / "No Motion" / "Inactive" entries (rows 171-183) ---
$pirSheet = $wb.Worksheets.Item("PIR")
$pirRows = @(
  @{r=171; A="2026-01-30"; B="17:10:12"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=172; A="2026-01-30"; B="17:10:12"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=173; A="2026-01-30"; B="17:10:17"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=174; A="2026-01-30"; B="17:10:22"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=175; A="2026-01-30"; B="17:10:27"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=176; A="2026-01-30"; B="17:10:32"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=177; A="2026-01-30"; B="17:10:37"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=178; A="2026-01-30"; B="17:10:42"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=179; A="2026-01-30"; B="17:10:47"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=180; A="2026-01-30"; B="17:10:53"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=181; A="2026-01-30"; B="17:10:58"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=182; A="2026-01-30"; B="17:11:03"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"},
  @{r=183; A="2026-01-30"; B="17:11:08"; C="17:00"; D="Bathroom"; E="No Motion"; F="Inactive"}
)
foreach ($row in $pirRows) {
    $pirSheet.Cells.Item($row.r, 1).NumberFormat = "@"
    $pirSheet.Cells.Item($row.r, 1).Value = $row.A
    $pirSheet.Cells.Item($row.r, 2).Value = $row.B
    $pirSheet.Cells.Item($row.r, 3).Value = $row.C
    $pirSheet.Cells.Item($row.r, 4).Value = $row.D
    if ($row.E -like "*%*") {
        $pirSheet.Cells.Item($row.r, 5).NumberFormat = "@"
    }
    $pirSheet.Cells.Item($row.r, 5).Value = $row.E
    $pirSheet.Cells.Item($row.r, 6).Value = $row.F
}

# --- Humidity sheet: new "Bathroom" humidity-% entries (rows 115-124) ---
$humiditySheet = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
  @{r=115; A="2026-01-30"; B="17:10:18"; C="17:00"; D="Bathroom"; E="86.4%"; F="Active"},
  @{r=116; A="2026-01-30"; B="17:10:23"; C="17:00"; D="Bathroom"; E="87.3%"; F="Active"},
  @{r=117; A="2026-01-30"; B="17:10:28"; C="17:00"; D="Bathroom"; E="86.4%"; F="Active"},
  @{r=118; A="2026-01-30"; B="17:10:33"; C="17:00"; D="Bathroom"; E="87.3%"; F="Active"},
  @{r=119; A="2026-01-30"; B="17:10:38"; C="17:00"; D="Bathroom"; E="86.4%"; F="Active"},
  @{r=120; A="2026-01-30"; B="17:10:43"; C="17:00"; D="Bathroom"; E="87.4%"; F="Active"},
  @{r=121; A="2026-01-30"; B="17:10:53"; C="17:00"; D="Bathroom"; E="87.3%"; F="Active"},
  @{r=122; A="2026-01-30"; B="17:10:58"; C="17:00"; D="Bathroom"; E="86.4%"; F="Active"},
  @{r=123; A="2026-01-30"; B="17:11:03"; C="17:00"; D="Bathroom"; E="87.3%"; F="Active"},
  @{r=124; A="2026-01-30"; B="17:11:08"; C="17:00"; D="Bathroom"; E="86.4%"; F="Active"}
)
foreach ($row in $humidityRows) {
    $humiditySheet.Cells.Item($row.r, 1).NumberFormat = "@"
    $humiditySheet.Cells.Item($row.r, 1).Value = $row.A
    $humiditySheet.Cells.Item($row.r, 2).Value = $row.B
    $humiditySheet.Cells.Item($row.r, 3).Value = $row.C
    $humiditySheet.Cells.Item($row.r, 4).Value = $row.D
    if ($row.E -like "*%*") {
        $humiditySheet.Cells.Item($row.r, 5).NumberFormat = "@"
    }
    $humiditySheet.Cells.Item($row.r, 5).Value = $row.E
    $humiditySheet.Cells.Item($row.r, 6).Value = $row.F
}
